$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.681.86'
$ws.Range("E2").Value = '  -1.74%  '

$ws.Range("D3").Value = '1.802.96'
$ws.Range("E3").Value = '  -1.26%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.94'
$ws.Range("E5").Value = '  -1.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5941'
$ws.Range("E6").Value = '  -1.60%  '

$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2778'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06824'
$ws.Range("E9").Value = '  -3.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.35'
$ws.Range("E10").Value = '  -0.65%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '1.959.40'
$ws.Range("E11").Value = '  +7.16%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07511'
$ws.Range("E12").Value = '  -1.55%  '

$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6229'
$ws.Range("E14").Value = '  -0.68%  '

$ws.Range("D15").Value = '2.047.78'
$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009221'
$ws.Range("E16").Value = '  -7.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '75.64'
$ws.Range("E17").Value = '  -4.23%  '

$ws.Range("D18").Value = '28.655.70'
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.480'
$ws.Range("E19").Value = '  -6.18%  '

$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '211.18'
$ws.Range("E21").Value = '  -6.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.52'
$ws.Range("E22").Value = '  -1.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.833'
$ws.Range("E23").Value = '  -2.37%  '

$ws.Range("E24").Value = '  +0.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.97'
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.859'
$ws.Range("E26").Value = '  -2.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1270'
$ws.Range("E27").Value = '  -2.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.46'
$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.410'
$ws.Range("E29").Value = '  -4.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06158'
$ws.Range("E30").Value = '  -1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.425'
$ws.Range("E31").Value = '  -1.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.781'
$ws.Range("E32").Value = '  -1.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.742'
$ws.Range("E33").Value = '  -1.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.731'
$ws.Range("E34").Value = '  -0.73%  '

$ws.Range("E35").Value = '  -5.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6427'
$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.502'
$ws.Range("E37").Value = '  -1.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.716'
$ws.Range("E38").Value = '  -0.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.550'
$ws.Range("E39").Value = '  +0.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01692'
$ws.Range("E40").Value = '  -2.05%  '

$ws.Range("D41").Value = '1.147.44'
$ws.Range("E41").Value = '  -5.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8886'
$ws.Range("E42").Value = '  -1.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.17'
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").Value = '1.952.11'
$ws.Range("E45").Value = '  -2.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.54'
$ws.Range("E46").Value = '  -3.46%  '

$ws.Range("E47").Value = '  -3.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.594'
$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.362'
$ws.Range("E49").Value = '  -1.37%  '

$ws.Range("E50").Value = '  -0.84%  '

$ws.Range("E51").Value = '  -1.66%  '
